# Add 2022-Q4 data
#
# Before: 总计 | 2022-Q3 | 2022-Q2
# After:  总计 | 2022-Q4 | 2022-Q3 | 2022-Q2
#
# The "总计" (totals) sheet gets a new row for 2022-Q4 (inserted above the
# existing 2022-Q3 / 2022-Q2 rows), and a brand new "2022-Q4" worksheet is
# inserted right before the existing "2022-Q3" sheet, holding the per-fund
# holdings detail for the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (totals) summary sheet: insert the 2022-Q4 row and
#    shift 2022-Q3 / 2022-Q2 down by one, renumbering the index column.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Push the old row 3 (2022-Q2) down into row 4, copying formats so the
# index cell (A4) keeps the same style as the other index cells.
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(4, 1).PasteSpecial(-4122)  # xlPasteFormats
$total.Cells.Item(4, 2).Value = "2022-Q2"
$total.Cells.Item(4, 3).Value = 5
$total.Cells.Item(4, 4).Value = 1

# Push the old row 2 (2022-Q3) down into row 3.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q3"
$total.Cells.Item(3, 3).Value = 10
$total.Cells.Item(3, 4).Value = 0.67

# Write the new 2022-Q4 row into row 2.
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 8
$total.Cells.Item(2, 4).Value = 0.83

# ---------------------------------------------------------------------
# 2. Create the new "2022-Q4" worksheet by duplicating the existing
#    "2022-Q3" sheet (so it inherits identical headers / column styles),
#    inserting it right before "2022-Q3", then overwriting its data.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The Q3 template has 10 data rows (rows 2-11); Q4 only needs 8 (rows 2-9),
# so drop the two trailing rows.
$q4.Rows.Item(10).Delete()
$q4.Rows.Item(10).Delete()

# Columns B (fund code) and D:G (scale/position/weight/value) are
# numeric-looking text (leading zeros, fixed decimal places) that must be
# preserved verbatim, so mark them as Text before writing the strings.
$q4.Range("B2:B9").NumberFormat = "@"
$q4.Range("D2:G9").NumberFormat = "@"

$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(2, 2).Value = "001877"
$q4.Cells.Item(2, 3).Value = "宝盈国家安全沪港深股票A"
$q4.Cells.Item(2, 4).Value = "6.58"
$q4.Cells.Item(2, 5).Value = "94.32"
$q4.Cells.Item(2, 6).Value = "4.67"
$q4.Cells.Item(2, 7).Value = "0.3073"
$q4.Cells.Item(2, 8).Value = 4

$q4.Cells.Item(3, 1).Value = 1
$q4.Cells.Item(3, 2).Value = "013895"
$q4.Cells.Item(3, 3).Value = "宝盈成长精选混合A"
$q4.Cells.Item(3, 4).Value = "8.59"
$q4.Cells.Item(3, 5).Value = "94.68"
$q4.Cells.Item(3, 6).Value = "3.23"
$q4.Cells.Item(3, 7).Value = "0.2775"
$q4.Cells.Item(3, 8).Value = 5

$q4.Cells.Item(4, 1).Value = 2
$q4.Cells.Item(4, 2).Value = "001075"
$q4.Cells.Item(4, 3).Value = "宝盈转型动力灵活配置混合A"
$q4.Cells.Item(4, 4).Value = "4.42"
$q4.Cells.Item(4, 5).Value = "91.93"
$q4.Cells.Item(4, 6).Value = "3.05"
$q4.Cells.Item(4, 7).Value = "0.1348"
$q4.Cells.Item(4, 8).Value = 9

$q4.Cells.Item(5, 1).Value = 3
$q4.Cells.Item(5, 2).Value = "013896"
$q4.Cells.Item(5, 3).Value = "宝盈成长精选混合C"
$q4.Cells.Item(5, 4).Value = "3.06"
$q4.Cells.Item(5, 5).Value = "94.68"
$q4.Cells.Item(5, 6).Value = "3.23"
$q4.Cells.Item(5, 7).Value = "0.0988"
$q4.Cells.Item(5, 8).Value = 5

$q4.Cells.Item(6, 1).Value = 4
$q4.Cells.Item(6, 2).Value = "013613"
$q4.Cells.Item(6, 3).Value = "宝盈国家安全沪港深股票C"
$q4.Cells.Item(6, 4).Value = "0.13"
$q4.Cells.Item(6, 5).Value = "94.32"
$q4.Cells.Item(6, 6).Value = "4.67"
$q4.Cells.Item(6, 7).Value = "0.0061"
$q4.Cells.Item(6, 8).Value = 4

$q4.Cells.Item(7, 1).Value = 5
$q4.Cells.Item(7, 2).Value = "015389"
$q4.Cells.Item(7, 3).Value = "宝盈转型动力灵活配置混合C"
$q4.Cells.Item(7, 4).Value = "0.03"
$q4.Cells.Item(7, 5).Value = "91.93"
$q4.Cells.Item(7, 6).Value = "3.05"
$q4.Cells.Item(7, 7).Value = "0.0009"
$q4.Cells.Item(7, 8).Value = 9

$q4.Cells.Item(8, 1).Value = 6
$q4.Cells.Item(8, 2).Value = "519222"
$q4.Cells.Item(8, 3).Value = "海富通欣益灵活配置混合A"
$q4.Cells.Item(8, 4).Value = "0.25"
$q4.Cells.Item(8, 5).Value = "31.65"
$q4.Cells.Item(8, 6).Value = "0.18"
$q4.Cells.Item(8, 7).Value = "0.0004"
$q4.Cells.Item(8, 8).Value = 3

$q4.Cells.Item(9, 1).Value = 7
$q4.Cells.Item(9, 2).Value = "519221"
$q4.Cells.Item(9, 3).Value = "海富通欣益灵活配置混合C"
$q4.Cells.Item(9, 4).Value = "0.10"
$q4.Cells.Item(9, 5).Value = "31.65"
$q4.Cells.Item(9, 6).Value = "0.18"
$q4.Cells.Item(9, 7).Value = "0.0002"
$q4.Cells.Item(9, 8).Value = 3

# ---------------------------------------------------------------------
# 3. Restore the original selection: "2022-Q2" (now the last tab) was the
#    active sheet before this edit.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$q2.Activate()
$q2.Range("A1").Select()
